$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: rent -> msa, 1000 -> 3000, date shifted
$ws.Range("A2").Value = "msa"
$ws.Range("B2").Value = 3000
$ws.Range("C2").Value = 45960.125185185185

# Row 3: additionaltest -> rent, amount stays 1000, date shifted
$ws.Range("A3").Value = "rent"
$ws.Range("B3").Value = 1000
$ws.Range("C3").Value = 45957.125185185185

# Row 4: additionaltest -> spris, 1000 -> 400, date shifted
$ws.Range("A4").Value = "spris"
$ws.Range("B4").Value = 400
$ws.Range("C4").Value = 45955.125185185185

# Row 5 (new): Room rent, 1000, same date as row4 -- copy format from C4 so it keeps the date style
$ws.Range("A5").Value = "Room rent"
$ws.Range("B5").Value = 1000
$ws.Range("C4").Copy()
$ws.Range("C5").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C5").Value = 45955.125185185185
